$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "233.07", "0.450") are not auto-converted to floating point
# numbers by the COM layer, preserving the original inline-string type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.941.49"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.252.04"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "233.07"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "0.643"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("D7").Value = "63.47"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.450"
$ws.Range("E9").Value = "  +5.70%  "
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("D11").Value = "57.87"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "26.52"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "2.586.23"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "15.58"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "6.10"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").Value = "0.835"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "2.251.00"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "43.818.94"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  +4.10%  "
$ws.Range("D21").Value = "73.17"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "6.10"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "249.61"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("E25").Value = "  -5.13%  "
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -7.30%  "
$ws.Range("D27").Value = "3.34"
$ws.Range("E27").Value = "  +21.45%  "
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "173.35"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "21.48"
$ws.Range("E30").Value = "  +4.35%  "
$ws.Range("D31").Value = "0.139"
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "1.43"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "0.126"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.82"
$ws.Range("E34").Value = "  +2.05%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0686"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("E36").Value = "  -5.02%  "
$ws.Range("D37").Value = "3.68"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "6.41"
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("E40").Value = "  +2.42%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "8.60"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "0.000224"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "97.92"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("D45").Value = "17.13"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "0.0949"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").Value = "4.37"
$ws.Range("E48").Value = "  -7.07%  "
$ws.Range("D49").Value = "1.442.84"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "9.84"
$ws.Range("E51").Value = "  -13.07%  "

# Restore the default (unstyled) cell style on column D now that the
# values have been written as text, matching the original workbook
# which had no explicit style override on these cells.
$ws.Range("D2:D51").Style = "Normal"
